# Update column C ("Förändrad") for rows 2-52: change date serial 45178 -> 45179
# (2023-09-09 -> 2023-09-10)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 52; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}
